# Insert a new price record as row 241 in the "Berenjena" price sheet.
# Excel's native row Insert() shifts the existing rows 241-273 down to
# 242-274 (carrying their values/styles with them), which matches the
# diff exactly. We only need to insert the row and fill in the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 241; everything below shifts down one.
$ws.Rows("241:241").Insert()

# Populate the newly inserted row 241 with the new weekly price record.
$ws.Range("A241").Value = 8
$ws.Range("B241").Value = "Terminal La Palmera de La Serena"
$ws.Range("C241").Value = "Coquimbo"
$ws.Range("D241").Value = 45131
$ws.Range("E241").Value = 4
$ws.Range("F241").Value = 100112001
$ws.Range("G241").Value = "Berenjena"
$ws.Range("H241").Value = "Sin especificar"
$ws.Range("I241").Value = "Primera"
$ws.Range("J241").Value = 420
$ws.Range("K241").Value = 8000
$ws.Range("L241").Value = 9000
$ws.Range("M241").Value = 8500
$ws.Range("N241").Value = "$/caja 50 unidades"
$ws.Range("O241").Value = "Región de Arica y Parinacota"
$ws.Range("P241").Value = 170
$ws.Range("Q241").Value = 50
$ws.Range("R241").Value = "Hortaliza"
